# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the Leve profit sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (45 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 481.5
$ws.Range("I2").Value = 481.5
$ws.Range("K2").Value = 481.5
$ws.Range("M2").Value = -368.5
$ws.Range("H17").Value = 1847.4762
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1847.4762
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5542.4286
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = -5878.4286
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null
$ws.Range("H106").Value = 52873.65
$ws.Range("J106").Value = 252641.5
$ws.Range("L106").Value = 252641.5
$ws.Range("N106").Value = -253903.5
$ws.Range("H107").Value = 736.75
$ws.Range("I107").Value = 704.5833
$ws.Range("K107").Value = 704.5833
$ws.Range("M107").Value = 1215.4167
$ws.Range("H111").Value = 17359.428
$ws.Range("I111").Value = 1497.8
$ws.Range("J111").Value = 57013.5
$ws.Range("K111").Value = 4493.4
$ws.Range("L111").Value = 171040.5
$ws.Range("M111").Value = -1426.4
$ws.Range("N111").Value = -177174.5
$ws.Range("H137").Value = 3628.6897
$ws.Range("J137").Value = 7368.8887
$ws.Range("L137").Value = 22106.6661
$ws.Range("N137").Value = -27206.6661
$ws.Range("H138").Value = 3307.6064
$ws.Range("I138").Value = 3179.5
$ws.Range("J138").Value = 3379.8718
$ws.Range("K138").Value = 9538.5
$ws.Range("L138").Value = 10139.6154
$ws.Range("M138").Value = -4398.5
$ws.Range("N138").Value = -20419.6154

# --- Sheet: ARM (18 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 700
$ws.Range("I12").Value = 700
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 700
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -527
$ws.Range("N12").Value = $null
$ws.Range("H32").Value = 18850.775
$ws.Range("I32").Value = 13444.7705
$ws.Range("J32").Value = 44799.6
$ws.Range("K32").Value = 13444.7705
$ws.Range("L32").Value = 44799.6
$ws.Range("M32").Value = -13157.7705
$ws.Range("N32").Value = -45373.6
$ws.Range("H111").Value = 61333.332
$ws.Range("J111").Value = 61333.332
$ws.Range("L111").Value = 61333.332
$ws.Range("N111").Value = -69513.33199999999

# --- Sheet: CRP (19 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 750
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null
$ws.Range("H31").Value = 4155.48
$ws.Range("I31").Value = 3045.2727
$ws.Range("K31").Value = 3045.2727
$ws.Range("M31").Value = -2750.2727
$ws.Range("H34").Value = 4155.48
$ws.Range("I34").Value = 3045.2727
$ws.Range("K34").Value = 3045.2727
$ws.Range("M34").Value = -2843.2727
$ws.Range("H132").Value = 215424.39
$ws.Range("I132").Value = 2441.5952
$ws.Range("J132").Value = 2004479.8
$ws.Range("K132").Value = 7324.785600000001
$ws.Range("L132").Value = 6013439.4
$ws.Range("M132").Value = -4794.785600000001
$ws.Range("N132").Value = -6018499.4

# --- Sheet: CUL (41 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16886150
$ws.Range("I4").Value = 335691.38
$ws.Range("K4").Value = 1007074.14
$ws.Range("M4").Value = -1006962.14
$ws.Range("H5").Value = 542.7895
$ws.Range("I5").Value = 522.94446
$ws.Range("K5").Value = 1568.83338
$ws.Range("M5").Value = -1456.83338
$ws.Range("H39").Value = 4422.273
$ws.Range("J39").Value = 4422.273
$ws.Range("L39").Value = 13266.819
$ws.Range("N39").Value = -13854.819
$ws.Range("H44").Value = 3250
$ws.Range("I44").Value = 3800
$ws.Range("J44").Value = 225
$ws.Range("K44").Value = 11400
$ws.Range("L44").Value = 675
$ws.Range("M44").Value = -11002
$ws.Range("N44").Value = -1471
$ws.Range("H112").Value = 7762.5
$ws.Range("I112").Value = 3400
$ws.Range("J112").Value = 9216.666999999999
$ws.Range("K112").Value = 10200
$ws.Range("L112").Value = 27650.001
$ws.Range("M112").Value = -9092
$ws.Range("N112").Value = -29866.001
$ws.Range("H114").Value = 10257
$ws.Range("I114").Value = 5514
$ws.Range("J114").Value = 15000
$ws.Range("K114").Value = 16542
$ws.Range("L114").Value = 45000
$ws.Range("M114").Value = -13288
$ws.Range("N114").Value = -51508
$ws.Range("H129").Value = 835.15
$ws.Range("J129").Value = 1441.1428
$ws.Range("L129").Value = 4323.428400000001
$ws.Range("N129").Value = -14323.4284
$ws.Range("H135").Value = 542.7895
$ws.Range("I135").Value = 522.94446
$ws.Range("K135").Value = 4706.50014
$ws.Range("M135").Value = -2171.50014

# --- Sheet: GSM (12 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 84997.5
$ws.Range("J63").Value = 84997.5
$ws.Range("L63").Value = 84997.5
$ws.Range("N63").Value = -86369.5
$ws.Range("H66").Value = 84997.5
$ws.Range("J66").Value = 84997.5
$ws.Range("L66").Value = 254992.5
$ws.Range("N66").Value = -261856.5
$ws.Range("H126").Value = 7576.1665
$ws.Range("I126").Value = 10350.8
$ws.Range("K126").Value = 31052.4
$ws.Range("M126").Value = -28582.4

# --- Sheet: LTW (12 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 84037.5
$ws.Range("J116").Value = 84037.5
$ws.Range("L116").Value = 84037.5
$ws.Range("N116").Value = -93215.5
$ws.Range("H132").Value = 119581.68
$ws.Range("I132").Value = 207440.75
$ws.Range("K132").Value = 622322.25
$ws.Range("M132").Value = -619792.25
$ws.Range("H136").Value = 7003.8945
$ws.Range("I136").Value = 7004.4707
$ws.Range("K136").Value = 21013.4121
$ws.Range("M136").Value = -18463.4121

# --- Sheet: WVR (22 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16673843
$ws.Range("I81").Value = 5353.6665
$ws.Range("J81").Value = 33342334
$ws.Range("K81").Value = 10707.333
$ws.Range("L81").Value = 66684668
$ws.Range("M81").Value = -9646.333000000001
$ws.Range("N81").Value = -66686790
$ws.Range("H84").Value = 16673843
$ws.Range("I84").Value = 5353.6665
$ws.Range("J84").Value = 33342334
$ws.Range("K84").Value = 53536.665
$ws.Range("L84").Value = 333423340
$ws.Range("M84").Value = -48232.665
$ws.Range("N84").Value = -333433948
$ws.Range("H132").Value = 776419.9
$ws.Range("I132").Value = 1119649.4
$ws.Range("K132").Value = 3358948.2
$ws.Range("M132").Value = -3356418.2
$ws.Range("H136").Value = 14042.526
$ws.Range("I136").Value = 18317.416
$ws.Range("K136").Value = 54952.24800000001
$ws.Range("M136").Value = -52402.24800000001

